$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = "68.245.20"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.642.51"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = "595.17"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = "159.62"
$ws.Range("E6").Value = "  +3.44%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.74%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = "5.26"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = "27.89"
$ws.Range("E13").Value = "  -0.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = "3.125.04"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  -2.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = "68.100.53"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = "2.627.60"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = "11.35"
$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = "359.57"
$ws.Range("E19").Value = "  -1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = "4.39"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = "4.73"
$ws.Range("E22").Value = "  -3.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = "75.05"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = "9.75"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value = "2.810.00"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D29").Value = "557.15"
$ws.Range("E29").Value = "  -2.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "7.98"
$ws.Range("E30").Value = "  -2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.37"
$ws.Range("E31").Value = "  -3.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.86"
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.127"
$ws.Range("E34").Value = "  -3.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -2.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "159.35"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D37").Value = "19.69"
$ws.Range("E37").Value = "  +1.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.370"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.85"
$ws.Range("E39").Value = "  -2.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "5.32"
$ws.Range("E40").Value = "  -1.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  -2.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").Value = "0.0₆0321"
$ws.Range("E42").Value = "  -5.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "156.61"
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.76"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "22.01"
$ws.Range("E46").Value = "  +0.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("B47").Value = "Optimism"
$ws.Range("C47").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D47").Value = "1.68"
$ws.Range("E47").Value = "  -2.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0773"
$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.612"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "0.565"
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("B51").Value = "Fantom"
$ws.Range("C51").Value = "https://coinranking.com/coin/uIEWfMFnQo9K_+fantom-ftm"
$ws.Range("D51").Value = "0.716"
$ws.Range("E51").Value = "  -2.43%  "
